# Apply the "Incidence of HT is (nearly) correct" revision:
#   - prevalence2018: replace the hard-coded probability constants in
#     C20:C82 with formulas that knock the raw STEPS-survey rate down by
#     an age-band-specific adherence/coverage factor.
#   - incidence2018_plus: the C20:C82 cells already derive from
#     prevalence2018!C.. / <divisor>; widen the divisor for the three
#     older age bands so the derived incidence keeps pace with the new
#     prevalence formulas.
#   - nudge the saved selection/scroll position on both sheets, matching
#     where the author's cursor ended up.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. prevalence2018 sheet: C20:C82
# ---------------------------------------------------------------------
$prev = $wb.Worksheets.Item("prevalence2018")

# Age 18-34 (rows 20-36): 4.9% raw rate, 92% coverage factor
$prev.Range("C20").Formula = "=0.049*0.92"
$prev.Range("C21:C36").Formula = "=0.049*0.92"

# Age 35-44 (rows 37-46): 6.2% raw rate, 94% coverage factor
$prev.Range("C37").Formula = "=0.062*0.94"
$prev.Range("C38:C46").Formula = "=0.062*0.94"

# Age 45-54 (rows 47-56): 5.6% raw rate, 87% coverage factor
$prev.Range("C47").Formula = "=0.056*0.87"
$prev.Range("C48:C56").Formula = "=0.056*0.87"

# Age 55-80 (rows 57-82): 6.8% raw rate, 89% coverage factor
$prev.Range("C57").Formula = "=0.068*0.89"
$prev.Range("C58:C82").Formula = "=0.068*0.89"

# ---------------------------------------------------------------------
# 2. incidence2018_plus sheet: widen the denominator for rows 37-82
# ---------------------------------------------------------------------
$inc = $wb.Worksheets.Item("incidence2018_plus")

for ($r = 37; $r -le 46; $r++) {
    $inc.Range("C$r").Formula = "=prevalence2018!C$r/11"
}

for ($r = 47; $r -le 56; $r++) {
    $inc.Range("C$r").Formula = "=prevalence2018!C$r/120"
}

for ($r = 57; $r -le 82; $r++) {
    $inc.Range("C$r").Formula = "=prevalence2018!C$r/200"
}

# ---------------------------------------------------------------------
# 3. Restore cursor / selection positions as left by the author
# ---------------------------------------------------------------------
$prev.Activate()
$prev.Range("F37").Select()

$inc.Activate()
$inc.Range("F53").Select()
